# Regenerate orders with updated distance/size labels.
#
# The stimulus set was relabeled: the three "Distance" conditions (D64,
# D80, D51) and the largest "Size" condition (S30) were renamed to
# D69, D86, D55 and S31 respectively. These tokens appear, as substrings,
# throughout the Condition / Filename_Left / Filename_Right / Distance /
# Size columns (e.g. "Face11_D64_S25" -> "Face11_D69_S25",
# "Fixation_D80_l.png" -> "Fixation_D86_l.png", the bare Distance value
# "D51" -> "D55", the bare Size value "S30" -> "S31"), so a simple global
# substring replace across the sheet's used range reproduces the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$used.Replace("D64", "D69")
$used.Replace("D80", "D86")
$used.Replace("D51", "D55")
$used.Replace("S30", "S31")
